$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 693.7778
$ws.Range("I5").Value = 693.7778
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 693.7778
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -578.7778
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 352.7
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5823.8237
$ws.Range("I62").Value = 2666.6667
$ws.Range("J62").Value = 9375.625
$ws.Range("K62").Value = 2666.6667
$ws.Range("L62").Value = 9375.625
$ws.Range("M62").Value = -2042.6667
$ws.Range("N62").Value = -10623.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5823.8237
$ws.Range("I65").Value = 2666.6667
$ws.Range("J65").Value = 9375.625
$ws.Range("K65").Value = 13333.3335
$ws.Range("L65").Value = 46878.125
$ws.Range("M65").Value = -10213.3335
$ws.Range("N65").Value = -53118.125

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3021.4412
$ws.Range("I86").Value = 2848.1765
$ws.Range("J86").Value = 3194.7058
$ws.Range("K86").Value = 2848.1765
$ws.Range("L86").Value = 3194.7058
$ws.Range("M86").Value = -1725.1765
$ws.Range("N86").Value = -5440.7058

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3021.4412
$ws.Range("I89").Value = 2848.1765
$ws.Range("J89").Value = 3194.7058
$ws.Range("K89").Value = 14240.8825
$ws.Range("L89").Value = 15973.529
$ws.Range("M89").Value = -8624.8825
$ws.Range("N89").Value = -27205.529

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2201.6
$ws.Range("I92").Value = 1076.6086
$ws.Range("J92").Value = 5898
$ws.Range("K92").Value = 1076.6086
$ws.Range("L92").Value = 5898
$ws.Range("M92").Value = 171.3914
$ws.Range("N92").Value = -8394

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1689
$ws.Range("J97").Value = 1689
$ws.Range("L97").Value = 5067
$ws.Range("N97").Value = -6059

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4744.9414
$ws.Range("I106").Value = 3332.4546
$ws.Range("K106").Value = 3332.4546
$ws.Range("M106").Value = -2701.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3679.8125
$ws.Range("I2").Value = 336.3846
$ws.Range("K2").Value = 336.3846
$ws.Range("M2").Value = -223.3846

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6132.8613
$ws.Range("I32").Value = 4734.6206
$ws.Range("K32").Value = 4734.6206
$ws.Range("M32").Value = -4447.6206

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3679.8125
$ws.Range("I116").Value = 336.3846
$ws.Range("K116").Value = 336.3846
$ws.Range("M116").Value = 1957.6154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3758.2778
$ws.Range("J122").Value = 4272.727
$ws.Range("L122").Value = 12818.181
$ws.Range("N122").Value = -17718.181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3679.8125
$ws.Range("I3").Value = 336.3846
$ws.Range("K3").Value = 336.3846
$ws.Range("M3").Value = -222.3846

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 3520.2
$ws.Range("I7").Value = 3900.25
$ws.Range("K7").Value = 3900.25
$ws.Range("M7").Value = -3787.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 10000
$ws.Range("J16").Value = 10000
$ws.Range("L16").Value = 10000
$ws.Range("N16").Value = -10340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2092.1428
$ws.Range("I22").Value = 432
$ws.Range("K22").Value = 432
$ws.Range("M22").Value = -82

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3551.6667
$ws.Range("I99").Value = 3404.8
$ws.Range("K99").Value = 3404.8
$ws.Range("M99").Value = -1906.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 884.1622
$ws.Range("I107").Value = 832.3913
$ws.Range("J107").Value = 969.2143
$ws.Range("K107").Value = 832.3913
$ws.Range("L107").Value = 969.2143
$ws.Range("M107").Value = 1087.6087
$ws.Range("N107").Value = -4809.2143

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3551.6667
$ws.Range("I126").Value = 3404.8
$ws.Range("K126").Value = 10214.4
$ws.Range("M126").Value = -7744.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2622.7188
$ws.Range("I132").Value = 1825.8077
$ws.Range("K132").Value = 5477.4231
$ws.Range("M132").Value = -2947.4231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 209.25
$ws.Range("I23").Value = 85
$ws.Range("J23").Value = 283.8
$ws.Range("K23").Value = 255
$ws.Range("L23").Value = 851.4000000000001
$ws.Range("M23").Value = -20
$ws.Range("N23").Value = -1321.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 23999
$ws.Range("I87").Value = 23999
$ws.Range("K87").Value = 71997
$ws.Range("M87").Value = -70749

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 23999
$ws.Range("I90").Value = 23999
$ws.Range("K90").Value = 215991
$ws.Range("M90").Value = -209751

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1571
$ws.Range("I113").Value = 1401.5
$ws.Range("J113").Value = 1643.6428
$ws.Range("K113").Value = 4204.5
$ws.Range("L113").Value = 4930.928400000001
$ws.Range("M113").Value = -2034.5
$ws.Range("N113").Value = -9270.928400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 2853
$ws.Range("I124").Value = 1680
$ws.Range("J124").Value = 3244
$ws.Range("K124").Value = 5040
$ws.Range("L124").Value = 9732
$ws.Range("M124").Value = -130
$ws.Range("N124").Value = -19552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1791.3077
$ws.Range("I140").Value = 1791.3077
$ws.Range("K140").Value = 5373.9231
$ws.Range("M140").Value = -193.9231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 9805.944
$ws.Range("J141").Value = 12321.818
$ws.Range("L141").Value = 36965.454
$ws.Range("N141").Value = -47325.454

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14294.238
$ws.Range("I70").Value = 10906.444
$ws.Range("J70").Value = 16835.084
$ws.Range("K70").Value = 10906.444
$ws.Range("L70").Value = 16835.084
$ws.Range("M70").Value = -10636.444
$ws.Range("N70").Value = -17375.084

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 14294.238
$ws.Range("I73").Value = 10906.444
$ws.Range("J73").Value = 16835.084
$ws.Range("K73").Value = 10906.444
$ws.Range("L73").Value = 16835.084
$ws.Range("M73").Value = -9970.444
$ws.Range("N73").Value = -18707.084

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2430.4243
$ws.Range("I102").Value = 1701.1428
$ws.Range("K102").Value = 1701.1428
$ws.Range("M102").Value = -79.14280000000008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 7079.143
$ws.Range("I132").Value = 6534.15
$ws.Range("K132").Value = 19602.45
$ws.Range("M132").Value = -17072.45

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8246.436
$ws.Range("I7").Value = 6632.52
$ws.Range("K7").Value = 6632.52
$ws.Range("M7").Value = -6520.52

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4959.577
$ws.Range("I22").Value = 2467.875
$ws.Range("J22").Value = 8946.299999999999
$ws.Range("K22").Value = 2467.875
$ws.Range("L22").Value = 8946.299999999999
$ws.Range("M22").Value = -2172.875
$ws.Range("N22").Value = -9536.299999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 4959.577
$ws.Range("I27").Value = 2467.875
$ws.Range("J27").Value = 8946.299999999999
$ws.Range("K27").Value = 2467.875
$ws.Range("L27").Value = 8946.299999999999
$ws.Range("M27").Value = -2360.875
$ws.Range("N27").Value = -9160.299999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7061.915
$ws.Range("I40").Value = 6828.0835
$ws.Range("K40").Value = 6828.0835
$ws.Range("M40").Value = -6692.0835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 128632.375
$ws.Range("I122").Value = 168846.42
$ws.Range("K122").Value = 506539.26
$ws.Range("M122").Value = -504089.26

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 8246.436
$ws.Range("I126").Value = 6632.52
$ws.Range("K126").Value = 19897.56
$ws.Range("M126").Value = -17427.56

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3656.8635
$ws.Range("I136").Value = 2926
$ws.Range("K136").Value = 8778
$ws.Range("M136").Value = -6228
